# Update metrics table (rows 2-26, columns A-Q) to reflect the retrained
# ensemble results ("atualizado todo o treinamento para o novo lm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 25,17

$arr[0,0] = "model_20_3_0"
$arr[0,1] = 0.9999805286793628
$arr[0,2] = 0.9991182490693772
$arr[0,3] = 0.9999115767816282
$arr[0,4] = 0.9998113399369803
$arr[0,5] = 0.9999463213469508
$arr[0,6] = 0.00001817562312264272
$arr[0,7] = 0.0008230757893437537
$arr[0,8] = 0.00005173478630590415
$arr[0,9] = 0.00006186285667827647
$arr[0,10] = 0.00005679882149209031
$arr[0,11] = 0.0002731861584103965
$arr[0,12] = 0.004263287830142685
$arr[0,13] = 1.000035947053484
$arr[0,14] = 0.004444784819294298
$arr[0,15] = 95.83085850106065
$arr[0,16] = 140.9292640211841

$arr[1,0] = "model_20_3_1"
$arr[1,1] = 0.9999805300863241
$arr[1,2] = 0.9991182456681967
$arr[1,3] = 0.9999115847276363
$arr[1,4] = 0.9998113507908325
$arr[1,5] = 0.9999463252255553
$arr[1,6] = 0.00001817430978600174
$arr[1,7] = 0.0008230789641964215
$arr[1,8] = 0.00005173013724387332
$arr[1,9] = 0.00006185929762983577
$arr[1,10] = 0.00005679471743685454
$arr[1,11] = 0.0002731710331822712
$arr[1,12] = 0.004263133798744973
$arr[1,13] = 1.000035944456017
$arr[1,14] = 0.004444624230461097
$arr[1,15] = 95.8310030225534
$arr[1,16] = 140.9294085426768

$arr[2,0] = "model_20_3_22"
$arr[2,1] = 0.9999805300082542
$arr[2,2] = 0.9991182316315311
$arr[2,3] = 0.9999115669109291
$arr[2,4] = 0.9998113827614422
$arr[2,5] = 0.9999463252535031
$arr[2,6] = 0.00001817438266093276
$arr[2,7] = 0.0008230920668076507
$arr[2,8] = 0.00005174056146905466
$arr[2,9] = 0.00006184881425984146
$arr[2,10] = 0.00005679468786444806
$arr[2,11] = 0.0002731515751168351
$arr[2,12] = 0.004263142345844525
$arr[2,13] = 1.000035944600146
$arr[2,14] = 0.004444633141428372
$arr[2,15] = 95.83099500301532
$arr[2,16] = 140.9294005231387

$arr[3,0] = "model_20_3_21"
$arr[3,1] = 0.9999805300082542
$arr[3,2] = 0.9991182316315311
$arr[3,3] = 0.9999115669109291
$arr[3,4] = 0.9998113827614422
$arr[3,5] = 0.9999463252535031
$arr[3,6] = 0.00001817438266093276
$arr[3,7] = 0.0008230920668076507
$arr[3,8] = 0.00005174056146905466
$arr[3,9] = 0.00006184881425984146
$arr[3,10] = 0.00005679468786444806
$arr[3,11] = 0.0002731515751168351
$arr[3,12] = 0.004263142345844525
$arr[3,13] = 1.000035944600146
$arr[3,14] = 0.004444633141428372
$arr[3,15] = 95.83099500301532
$arr[3,16] = 140.9294005231387

$arr[4,0] = "model_20_3_20"
$arr[4,1] = 0.9999805300082542
$arr[4,2] = 0.9991182316315311
$arr[4,3] = 0.9999115669109291
$arr[4,4] = 0.9998113827614422
$arr[4,5] = 0.9999463252535031
$arr[4,6] = 0.00001817438266093276
$arr[4,7] = 0.0008230920668076507
$arr[4,8] = 0.00005174056146905466
$arr[4,9] = 0.00006184881425984146
$arr[4,10] = 0.00005679468786444806
$arr[4,11] = 0.0002731515751168351
$arr[4,12] = 0.004263142345844525
$arr[4,13] = 1.000035944600146
$arr[4,14] = 0.004444633141428372
$arr[4,15] = 95.83099500301532
$arr[4,16] = 140.9294005231387

$arr[5,0] = "model_20_3_19"
$arr[5,1] = 0.9999805300082542
$arr[5,2] = 0.9991182316315311
$arr[5,3] = 0.9999115669109291
$arr[5,4] = 0.9998113827614422
$arr[5,5] = 0.9999463252535031
$arr[5,6] = 0.00001817438266093276
$arr[5,7] = 0.0008230920668076507
$arr[5,8] = 0.00005174056146905466
$arr[5,9] = 0.00006184881425984146
$arr[5,10] = 0.00005679468786444806
$arr[5,11] = 0.0002731515751168351
$arr[5,12] = 0.004263142345844525
$arr[5,13] = 1.000035944600146
$arr[5,14] = 0.004444633141428372
$arr[5,15] = 95.83099500301532
$arr[5,16] = 140.9294005231387

$arr[6,0] = "model_20_3_18"
$arr[6,1] = 0.9999805300082542
$arr[6,2] = 0.9991182316315311
$arr[6,3] = 0.9999115669109291
$arr[6,4] = 0.9998113827614422
$arr[6,5] = 0.9999463252535031
$arr[6,6] = 0.00001817438266093276
$arr[6,7] = 0.0008230920668076507
$arr[6,8] = 0.00005174056146905466
$arr[6,9] = 0.00006184881425984146
$arr[6,10] = 0.00005679468786444806
$arr[6,11] = 0.0002731515751168351
$arr[6,12] = 0.004263142345844525
$arr[6,13] = 1.000035944600146
$arr[6,14] = 0.004444633141428372
$arr[6,15] = 95.83099500301532
$arr[6,16] = 140.9294005231387

$arr[7,0] = "model_20_3_17"
$arr[7,1] = 0.9999805300082542
$arr[7,2] = 0.9991182316315311
$arr[7,3] = 0.9999115669109291
$arr[7,4] = 0.9998113827614422
$arr[7,5] = 0.9999463252535031
$arr[7,6] = 0.00001817438266093276
$arr[7,7] = 0.0008230920668076507
$arr[7,8] = 0.00005174056146905466
$arr[7,9] = 0.00006184881425984146
$arr[7,10] = 0.00005679468786444806
$arr[7,11] = 0.0002731515751168351
$arr[7,12] = 0.004263142345844525
$arr[7,13] = 1.000035944600146
$arr[7,14] = 0.004444633141428372
$arr[7,15] = 95.83099500301532
$arr[7,16] = 140.9294005231387

$arr[8,0] = "model_20_3_16"
$arr[8,1] = 0.9999805300082542
$arr[8,2] = 0.9991182316315311
$arr[8,3] = 0.9999115669109291
$arr[8,4] = 0.9998113827614422
$arr[8,5] = 0.9999463252535031
$arr[8,6] = 0.00001817438266093276
$arr[8,7] = 0.0008230920668076507
$arr[8,8] = 0.00005174056146905466
$arr[8,9] = 0.00006184881425984146
$arr[8,10] = 0.00005679468786444806
$arr[8,11] = 0.0002731515751168351
$arr[8,12] = 0.004263142345844525
$arr[8,13] = 1.000035944600146
$arr[8,14] = 0.004444633141428372
$arr[8,15] = 95.83099500301532
$arr[8,16] = 140.9294005231387

$arr[9,0] = "model_20_3_15"
$arr[9,1] = 0.9999805300082542
$arr[9,2] = 0.9991182316315311
$arr[9,3] = 0.9999115669109291
$arr[9,4] = 0.9998113827614422
$arr[9,5] = 0.9999463252535031
$arr[9,6] = 0.00001817438266093276
$arr[9,7] = 0.0008230920668076507
$arr[9,8] = 0.00005174056146905466
$arr[9,9] = 0.00006184881425984146
$arr[9,10] = 0.00005679468786444806
$arr[9,11] = 0.0002731515751168351
$arr[9,12] = 0.004263142345844525
$arr[9,13] = 1.000035944600146
$arr[9,14] = 0.004444633141428372
$arr[9,15] = 95.83099500301532
$arr[9,16] = 140.9294005231387

$arr[10,0] = "model_20_3_14"
$arr[10,1] = 0.9999805300082542
$arr[10,2] = 0.9991182316315311
$arr[10,3] = 0.9999115669109291
$arr[10,4] = 0.9998113827614422
$arr[10,5] = 0.9999463252535031
$arr[10,6] = 0.00001817438266093276
$arr[10,7] = 0.0008230920668076507
$arr[10,8] = 0.00005174056146905466
$arr[10,9] = 0.00006184881425984146
$arr[10,10] = 0.00005679468786444806
$arr[10,11] = 0.0002731515751168351
$arr[10,12] = 0.004263142345844525
$arr[10,13] = 1.000035944600146
$arr[10,14] = 0.004444633141428372
$arr[10,15] = 95.83099500301532
$arr[10,16] = 140.9294005231387

$arr[11,0] = "model_20_3_13"
$arr[11,1] = 0.9999805300082542
$arr[11,2] = 0.9991182316315311
$arr[11,3] = 0.9999115669109291
$arr[11,4] = 0.9998113827614422
$arr[11,5] = 0.9999463252535031
$arr[11,6] = 0.00001817438266093276
$arr[11,7] = 0.0008230920668076507
$arr[11,8] = 0.00005174056146905466
$arr[11,9] = 0.00006184881425984146
$arr[11,10] = 0.00005679468786444806
$arr[11,11] = 0.0002731515751168351
$arr[11,12] = 0.004263142345844525
$arr[11,13] = 1.000035944600146
$arr[11,14] = 0.004444633141428372
$arr[11,15] = 95.83099500301532
$arr[11,16] = 140.9294005231387

$arr[12,0] = "model_20_3_12"
$arr[12,1] = 0.9999805300082542
$arr[12,2] = 0.9991182316315311
$arr[12,3] = 0.9999115669109291
$arr[12,4] = 0.9998113827614422
$arr[12,5] = 0.9999463252535031
$arr[12,6] = 0.00001817438266093276
$arr[12,7] = 0.0008230920668076507
$arr[12,8] = 0.00005174056146905466
$arr[12,9] = 0.00006184881425984146
$arr[12,10] = 0.00005679468786444806
$arr[12,11] = 0.0002731515751168351
$arr[12,12] = 0.004263142345844525
$arr[12,13] = 1.000035944600146
$arr[12,14] = 0.004444633141428372
$arr[12,15] = 95.83099500301532
$arr[12,16] = 140.9294005231387

$arr[13,0] = "model_20_3_11"
$arr[13,1] = 0.9999805300082542
$arr[13,2] = 0.9991182316315311
$arr[13,3] = 0.9999115669109291
$arr[13,4] = 0.9998113827614422
$arr[13,5] = 0.9999463252535031
$arr[13,6] = 0.00001817438266093276
$arr[13,7] = 0.0008230920668076507
$arr[13,8] = 0.00005174056146905466
$arr[13,9] = 0.00006184881425984146
$arr[13,10] = 0.00005679468786444806
$arr[13,11] = 0.0002731515751168351
$arr[13,12] = 0.004263142345844525
$arr[13,13] = 1.000035944600146
$arr[13,14] = 0.004444633141428372
$arr[13,15] = 95.83099500301532
$arr[13,16] = 140.9294005231387

$arr[14,0] = "model_20_3_10"
$arr[14,1] = 0.9999805300082542
$arr[14,2] = 0.9991182316315311
$arr[14,3] = 0.9999115669109291
$arr[14,4] = 0.9998113827614422
$arr[14,5] = 0.9999463252535031
$arr[14,6] = 0.00001817438266093276
$arr[14,7] = 0.0008230920668076507
$arr[14,8] = 0.00005174056146905466
$arr[14,9] = 0.00006184881425984146
$arr[14,10] = 0.00005679468786444806
$arr[14,11] = 0.0002731515751168351
$arr[14,12] = 0.004263142345844525
$arr[14,13] = 1.000035944600146
$arr[14,14] = 0.004444633141428372
$arr[14,15] = 95.83099500301532
$arr[14,16] = 140.9294005231387

$arr[15,0] = "model_20_3_9"
$arr[15,1] = 0.9999805300082542
$arr[15,2] = 0.9991182316315311
$arr[15,3] = 0.9999115669109291
$arr[15,4] = 0.9998113827614422
$arr[15,5] = 0.9999463252535031
$arr[15,6] = 0.00001817438266093276
$arr[15,7] = 0.0008230920668076507
$arr[15,8] = 0.00005174056146905466
$arr[15,9] = 0.00006184881425984146
$arr[15,10] = 0.00005679468786444806
$arr[15,11] = 0.0002731515751168351
$arr[15,12] = 0.004263142345844525
$arr[15,13] = 1.000035944600146
$arr[15,14] = 0.004444633141428372
$arr[15,15] = 95.83099500301532
$arr[15,16] = 140.9294005231387

$arr[16,0] = "model_20_3_8"
$arr[16,1] = 0.9999805300082542
$arr[16,2] = 0.9991182316315311
$arr[16,3] = 0.9999115669109291
$arr[16,4] = 0.9998113827614422
$arr[16,5] = 0.9999463252535031
$arr[16,6] = 0.00001817438266093276
$arr[16,7] = 0.0008230920668076507
$arr[16,8] = 0.00005174056146905466
$arr[16,9] = 0.00006184881425984146
$arr[16,10] = 0.00005679468786444806
$arr[16,11] = 0.0002731515751168351
$arr[16,12] = 0.004263142345844525
$arr[16,13] = 1.000035944600146
$arr[16,14] = 0.004444633141428372
$arr[16,15] = 95.83099500301532
$arr[16,16] = 140.9294005231387

$arr[17,0] = "model_20_3_7"
$arr[17,1] = 0.9999805300082542
$arr[17,2] = 0.9991182316315311
$arr[17,3] = 0.9999115669109291
$arr[17,4] = 0.9998113827614422
$arr[17,5] = 0.9999463252535031
$arr[17,6] = 0.00001817438266093276
$arr[17,7] = 0.0008230920668076507
$arr[17,8] = 0.00005174056146905466
$arr[17,9] = 0.00006184881425984146
$arr[17,10] = 0.00005679468786444806
$arr[17,11] = 0.0002731515751168351
$arr[17,12] = 0.004263142345844525
$arr[17,13] = 1.000035944600146
$arr[17,14] = 0.004444633141428372
$arr[17,15] = 95.83099500301532
$arr[17,16] = 140.9294005231387

$arr[18,0] = "model_20_3_6"
$arr[18,1] = 0.9999805300082542
$arr[18,2] = 0.9991182316315311
$arr[18,3] = 0.9999115669109291
$arr[18,4] = 0.9998113827614422
$arr[18,5] = 0.9999463252535031
$arr[18,6] = 0.00001817438266093276
$arr[18,7] = 0.0008230920668076507
$arr[18,8] = 0.00005174056146905466
$arr[18,9] = 0.00006184881425984146
$arr[18,10] = 0.00005679468786444806
$arr[18,11] = 0.0002731515751168351
$arr[18,12] = 0.004263142345844525
$arr[18,13] = 1.000035944600146
$arr[18,14] = 0.004444633141428372
$arr[18,15] = 95.83099500301532
$arr[18,16] = 140.9294005231387

$arr[19,0] = "model_20_3_5"
$arr[19,1] = 0.9999805300082542
$arr[19,2] = 0.9991182316315311
$arr[19,3] = 0.9999115669109291
$arr[19,4] = 0.9998113827614422
$arr[19,5] = 0.9999463252535031
$arr[19,6] = 0.00001817438266093276
$arr[19,7] = 0.0008230920668076507
$arr[19,8] = 0.00005174056146905466
$arr[19,9] = 0.00006184881425984146
$arr[19,10] = 0.00005679468786444806
$arr[19,11] = 0.0002731515751168351
$arr[19,12] = 0.004263142345844525
$arr[19,13] = 1.000035944600146
$arr[19,14] = 0.004444633141428372
$arr[19,15] = 95.83099500301532
$arr[19,16] = 140.9294005231387

$arr[20,0] = "model_20_3_4"
$arr[20,1] = 0.9999805300082542
$arr[20,2] = 0.9991182316315311
$arr[20,3] = 0.9999115669109291
$arr[20,4] = 0.9998113827614422
$arr[20,5] = 0.9999463252535031
$arr[20,6] = 0.00001817438266093276
$arr[20,7] = 0.0008230920668076507
$arr[20,8] = 0.00005174056146905466
$arr[20,9] = 0.00006184881425984146
$arr[20,10] = 0.00005679468786444806
$arr[20,11] = 0.0002731515751168351
$arr[20,12] = 0.004263142345844525
$arr[20,13] = 1.000035944600146
$arr[20,14] = 0.004444633141428372
$arr[20,15] = 95.83099500301532
$arr[20,16] = 140.9294005231387

$arr[21,0] = "model_20_3_3"
$arr[21,1] = 0.9999805300082542
$arr[21,2] = 0.9991182316315311
$arr[21,3] = 0.9999115669109291
$arr[21,4] = 0.9998113827614422
$arr[21,5] = 0.9999463252535031
$arr[21,6] = 0.00001817438266093276
$arr[21,7] = 0.0008230920668076507
$arr[21,8] = 0.00005174056146905466
$arr[21,9] = 0.00006184881425984146
$arr[21,10] = 0.00005679468786444806
$arr[21,11] = 0.0002731515751168351
$arr[21,12] = 0.004263142345844525
$arr[21,13] = 1.000035944600146
$arr[21,14] = 0.004444633141428372
$arr[21,15] = 95.83099500301532
$arr[21,16] = 140.9294005231387

$arr[22,0] = "model_20_3_2"
$arr[22,1] = 0.9999805300082542
$arr[22,2] = 0.9991182316315311
$arr[22,3] = 0.9999115669109291
$arr[22,4] = 0.9998113827614422
$arr[22,5] = 0.9999463252535031
$arr[22,6] = 0.00001817438266093276
$arr[22,7] = 0.0008230920668076507
$arr[22,8] = 0.00005174056146905466
$arr[22,9] = 0.00006184881425984146
$arr[22,10] = 0.00005679468786444806
$arr[22,11] = 0.0002731515751168351
$arr[22,12] = 0.004263142345844525
$arr[22,13] = 1.000035944600146
$arr[22,14] = 0.004444633141428372
$arr[22,15] = 95.83099500301532
$arr[22,16] = 140.9294005231387

$arr[23,0] = "model_20_3_23"
$arr[23,1] = 0.9999805300082542
$arr[23,2] = 0.9991182316315311
$arr[23,3] = 0.9999115669109291
$arr[23,4] = 0.9998113827614422
$arr[23,5] = 0.9999463252535031
$arr[23,6] = 0.00001817438266093276
$arr[23,7] = 0.0008230920668076507
$arr[23,8] = 0.00005174056146905466
$arr[23,9] = 0.00006184881425984146
$arr[23,10] = 0.00005679468786444806
$arr[23,11] = 0.0002731515751168351
$arr[23,12] = 0.004263142345844525
$arr[23,13] = 1.000035944600146
$arr[23,14] = 0.004444633141428372
$arr[23,15] = 95.83099500301532
$arr[23,16] = 140.9294005231387

$arr[24,0] = "model_20_3_24"
$arr[24,1] = 0.9999805300082542
$arr[24,2] = 0.9991182316315311
$arr[24,3] = 0.9999115669109291
$arr[24,4] = 0.9998113827614422
$arr[24,5] = 0.9999463252535031
$arr[24,6] = 0.00001817438266093276
$arr[24,7] = 0.0008230920668076507
$arr[24,8] = 0.00005174056146905466
$arr[24,9] = 0.00006184881425984146
$arr[24,10] = 0.00005679468786444806
$arr[24,11] = 0.0002731515751168351
$arr[24,12] = 0.004263142345844525
$arr[24,13] = 1.000035944600146
$arr[24,14] = 0.004444633141428372
$arr[24,15] = 95.83099500301532
$arr[24,16] = 140.9294005231387

$ws.Range("A2:Q26").Value2 = $arr
